# Error Calculations and Plots
# Applies the missing-data imputation edits to Sheet1:
#   - removes the "RM 232" row and the "SC 92" row entirely
#   - fills in several previously-blank cells with numeric values
#   - blanks out several previously-numeric cells (swap of which values
#     are "missing" in this imputation combination)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the two rows that are removed outright.
# Delete row 26 ("RM 232") first; "SC 92" (originally row 28) becomes
# row 27 once row 26 is gone, so delete row 27 next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# 2) Set cells that go from blank -> numeric value.
$ws.Range("D2").Value = -13.5
$ws.Range("E4").Value = -6.4
$ws.Range("D12").Value = -14.1
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = -8.5
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("E23").Value = -7
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = -20.2
$ws.Range("B30").Value = -19.7
$ws.Range("D31").Value = -13.7
$ws.Range("D33").Value = -14.1

# 3) Clear cells that go from numeric value -> blank.
$ws.Range("E3").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("B27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("B32").ClearContents()

"Done"
